$wb = $excel.ActiveWorkbook

# sheet1 ("展览") - column F ("想去人数") updates
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 601  # F4: 599 -> 601
$ws.Cells.Item(5, 6).Value = 2617  # F5: 2610 -> 2617
$ws.Cells.Item(10, 6).Value = 5639  # F10: 5615 -> 5639
$ws.Cells.Item(11, 6).Value = 898  # F11: 899 -> 898
$ws.Cells.Item(13, 6).Value = 1491  # F13: 1488 -> 1491
$ws.Cells.Item(14, 6).Value = 1426  # F14: 1424 -> 1426
$ws.Cells.Item(15, 6).Value = 615  # F15: 614 -> 615
$ws.Cells.Item(17, 6).Value = 53  # F17: 52 -> 53
$ws.Cells.Item(19, 6).Value = 4758  # F19: 4751 -> 4758
$ws.Cells.Item(22, 6).Value = 2390  # F22: 2387 -> 2390
$ws.Cells.Item(23, 6).Value = 1275  # F23: 1273 -> 1275
$ws.Cells.Item(24, 6).Value = 459  # F24: 456 -> 459
$ws.Cells.Item(25, 6).Value = 1176  # F25: 1174 -> 1176
$ws.Cells.Item(26, 6).Value = 236  # F26: 234 -> 236
$ws.Cells.Item(27, 6).Value = 95  # F27: 94 -> 95
$ws.Cells.Item(28, 6).Value = 88  # F28: 85 -> 88
$ws.Cells.Item(29, 6).Value = 175  # F29: 173 -> 175
$ws.Cells.Item(30, 6).Value = 372  # F30: 371 -> 372
$ws.Cells.Item(31, 6).Value = 1296  # F31: 1294 -> 1296
$ws.Cells.Item(32, 6).Value = 2008  # F32: 2007 -> 2008
$ws.Cells.Item(35, 6).Value = 17  # F35: 15 -> 17
$ws.Cells.Item(37, 6).Value = 1390  # F37: 1388 -> 1390
$ws.Cells.Item(38, 6).Value = 596  # F38: 595 -> 596
$ws.Cells.Item(39, 6).Value = 93  # F39: 92 -> 93
$ws.Cells.Item(40, 6).Value = 528  # F40: 525 -> 528
$ws.Cells.Item(41, 6).Value = 189  # F41: 186 -> 189
$ws.Cells.Item(42, 6).Value = 1645  # F42: 1644 -> 1645
$ws.Cells.Item(43, 6).Value = 2434  # F43: 2432 -> 2434
$ws.Cells.Item(45, 6).Value = 80  # F45: 79 -> 80
$ws.Cells.Item(48, 6).Value = 31  # F48: 30 -> 31

# sheet2 ("演出") - column F ("想去人数") updates
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 7  # F11: 6 -> 7
$ws.Cells.Item(16, 6).Value = 188  # F16: 187 -> 188
$ws.Cells.Item(28, 6).Value = 296  # F28: 294 -> 296

# sheet3 ("本地生活") - column F ("想去人数") updates
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(6, 6).Value = 1670  # F6: 1668 -> 1670
$ws.Cells.Item(7, 6).Value = 542  # F7: 541 -> 542
$ws.Cells.Item(8, 6).Value = 1340  # F8: 1333 -> 1340
$ws.Cells.Item(9, 6).Value = 1201  # F9: 1202 -> 1201
$ws.Cells.Item(10, 6).Value = 1765  # F10: 1764 -> 1765
$ws.Cells.Item(11, 6).Value = 2255  # F11: 2246 -> 2255
$ws.Cells.Item(12, 6).Value = 693  # F12: 688 -> 693
$ws.Cells.Item(13, 6).Value = 567  # F13: 563 -> 567

# sheet4 ("全部类型") - column F ("想去人数") updates
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 1670  # F3: 1668 -> 1670
$ws.Cells.Item(5, 6).Value = 601  # F5: 599 -> 601
$ws.Cells.Item(6, 6).Value = 542  # F6: 541 -> 542
$ws.Cells.Item(7, 6).Value = 2617  # F7: 2610 -> 2617
$ws.Cells.Item(9, 6).Value = 1340  # F9: 1333 -> 1340
$ws.Cells.Item(11, 6).Value = 2255  # F11: 2246 -> 2255
$ws.Cells.Item(12, 6).Value = 5639  # F12: 5616 -> 5639
$ws.Cells.Item(13, 6).Value = 693  # F13: 688 -> 693
$ws.Cells.Item(17, 6).Value = 1491  # F17: 1488 -> 1491
$ws.Cells.Item(18, 6).Value = 1426  # F18: 1424 -> 1426
$ws.Cells.Item(20, 6).Value = 567  # F20: 563 -> 567
$ws.Cells.Item(21, 6).Value = 53  # F21: 52 -> 53
$ws.Cells.Item(22, 6).Value = 4758  # F22: 4751 -> 4758
$ws.Cells.Item(23, 6).Value = 2390  # F23: 2387 -> 2390
$ws.Cells.Item(24, 6).Value = 1275  # F24: 1273 -> 1275
$ws.Cells.Item(25, 6).Value = 459  # F25: 456 -> 459
$ws.Cells.Item(26, 6).Value = 1176  # F26: 1174 -> 1176
$ws.Cells.Item(27, 6).Value = 236  # F27: 234 -> 236
$ws.Cells.Item(28, 6).Value = 88  # F28: 85 -> 88
$ws.Cells.Item(30, 6).Value = 175  # F30: 173 -> 175
$ws.Cells.Item(32, 6).Value = 188  # F32: 187 -> 188
$ws.Cells.Item(33, 6).Value = 372  # F33: 371 -> 372
$ws.Cells.Item(34, 6).Value = 2008  # F34: 2007 -> 2008
$ws.Cells.Item(38, 6).Value = 1390  # F38: 1388 -> 1390
$ws.Cells.Item(39, 6).Value = 528  # F39: 525 -> 528
$ws.Cells.Item(42, 6).Value = 189  # F42: 186 -> 189
$ws.Cells.Item(44, 6).Value = 1645  # F44: 1644 -> 1645
$ws.Cells.Item(45, 6).Value = 2434  # F45: 2432 -> 2434
$ws.Cells.Item(46, 6).Value = 80  # F46: 79 -> 80
$ws.Cells.Item(49, 6).Value = 31  # F49: 30 -> 31
